# Auto-generated Excel COM-interop script
# Applies numeric corrections to currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the 8 job sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 80
$ws.Range("H80").Value = 2666
$ws.Range("I80").Value = 2242.875
$ws.Range("J80").Value = 3149.5715
$ws.Range("K80").Value = 6728.625
$ws.Range("L80").Value = 9448.7145
$ws.Range("M80").Value = -5730.625
$ws.Range("N80").Value = -11444.7145

# Row 83
$ws.Range("H83").Value = 2666
$ws.Range("I83").Value = 2242.875
$ws.Range("J83").Value = 3149.5715
$ws.Range("K83").Value = 20185.875
$ws.Range("L83").Value = 28346.1435
$ws.Range("M83").Value = -15193.875
$ws.Range("N83").Value = -38330.1435

# Row 86
$ws.Range("H86").Value = 14349.75
$ws.Range("I86").Value = 2699.5
$ws.Range("J86").Value = 26000
$ws.Range("K86").Value = 2699.5
$ws.Range("L86").Value = 26000
$ws.Range("M86").Value = -1576.5
$ws.Range("N86").Value = -28246

# Row 89
$ws.Range("H89").Value = 14349.75
$ws.Range("I89").Value = 2699.5
$ws.Range("J89").Value = 26000
$ws.Range("K89").Value = 13497.5
$ws.Range("L89").Value = 130000
$ws.Range("M89").Value = -7881.5
$ws.Range("N89").Value = -141232

# Row 132
$ws.Range("H132").Value = 10813.19
$ws.Range("I132").Value = 10813.19
$ws.Range("K132").Value = 32439.57
$ws.Range("M132").Value = -29909.57

# Row 138
$ws.Range("H138").Value = 2512.8628
$ws.Range("I138").Value = 2793.2856
$ws.Range("J138").Value = 2171.4783
$ws.Range("K138").Value = 8379.856800000001
$ws.Range("L138").Value = 6514.4349
$ws.Range("M138").Value = -3239.856800000001
$ws.Range("N138").Value = -16794.4349

$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 3049.6667
$ws.Range("I61").Value = 2766
$ws.Range("K61").Value = 2766
$ws.Range("M61").Value = -2554

# Row 74
$ws.Range("H74").Value = 5848.314
$ws.Range("I74").Value = 3782
$ws.Range("K74").Value = 3782
$ws.Range("M74").Value = -2908

# Row 77
$ws.Range("H77").Value = 5848.314
$ws.Range("I77").Value = 3782
$ws.Range("K77").Value = 18910
$ws.Range("M77").Value = -14542

# Row 132
$ws.Range("H132").Value = 4722.616
$ws.Range("I132").Value = 3189.8125
$ws.Range("K132").Value = 9569.4375
$ws.Range("M132").Value = -7039.4375

# Row 136
$ws.Range("H136").Value = 3049.6667
$ws.Range("I136").Value = 2766
$ws.Range("K136").Value = 8298
$ws.Range("M136").Value = -5748

$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 1951.6111
$ws.Range("I86").Value = 1804.9375
$ws.Range("J86").Value = 3125
$ws.Range("K86").Value = 1804.9375
$ws.Range("L86").Value = 3125
$ws.Range("M86").Value = -681.9375
$ws.Range("N86").Value = -5371

# Row 89
$ws.Range("H89").Value = 1951.6111
$ws.Range("I89").Value = 1804.9375
$ws.Range("J89").Value = 3125
$ws.Range("K89").Value = 9024.6875
$ws.Range("L89").Value = 15625
$ws.Range("M89").Value = -3408.6875
$ws.Range("N89").Value = -26857

# Row 94
$ws.Range("H94").Value = 10620
$ws.Range("I94").Value = 6033.3335
$ws.Range("K94").Value = 6033.3335
$ws.Range("M94").Value = -5582.3335

# Row 107
$ws.Range("H107").Value = 1200
$ws.Range("I107").Value = 1200
$ws.Range("K107").Value = 1200
$ws.Range("M107").Value = 720

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 134
$ws.Range("H134").Value = 3233.4167
$ws.Range("I134").Value = 3233.4167
$ws.Range("K134").Value = 9700.250100000001
$ws.Range("M134").Value = -7165.250100000001

$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2771.6365
$ws.Range("I31").Value = 3061
$ws.Range("K31").Value = 3061
$ws.Range("M31").Value = -2766

# Row 34
$ws.Range("H34").Value = 2771.6365
$ws.Range("I34").Value = 3061
$ws.Range("K34").Value = 3061
$ws.Range("M34").Value = -2859

# Row 58
$ws.Range("H58").Value = 3297.6453
$ws.Range("I58").Value = 3753.0667
$ws.Range("J58").Value = 2870.6875
$ws.Range("K58").Value = 3753.0667
$ws.Range("L58").Value = 2870.6875
$ws.Range("M58").Value = -3550.0667
$ws.Range("N58").Value = -3276.6875

# Row 62
$ws.Range("H62").Value = 2283.3333
$ws.Range("I62").Value = 1900
$ws.Range("K62").Value = 1900
$ws.Range("M62").Value = -1276

# Row 65
$ws.Range("H65").Value = 2283.3333
$ws.Range("I65").Value = 1900
$ws.Range("K65").Value = 9500
$ws.Range("M65").Value = -6380

# Row 107
$ws.Range("H107").Value = 846.875
$ws.Range("I107").Value = 635.0714
$ws.Range("K107").Value = 635.0714
$ws.Range("M107").Value = 1284.9286

# Row 132
$ws.Range("H132").Value = 1945.7678
$ws.Range("I132").Value = 1762.8077
$ws.Range("K132").Value = 5288.4231
$ws.Range("M132").Value = -2758.4231

# Row 134
$ws.Range("H134").Value = 1537.2609
$ws.Range("I134").Value = 1555.3684
$ws.Range("J134").Value = 1451.25
$ws.Range("K134").Value = 4666.1052
$ws.Range("L134").Value = 4353.75
$ws.Range("M134").Value = -2131.1052
$ws.Range("N134").Value = -9423.75

# Row 136
$ws.Range("H136").Value = 3297.6453
$ws.Range("I136").Value = 3753.0667
$ws.Range("J136").Value = 2870.6875
$ws.Range("K136").Value = 11259.2001
$ws.Range("L136").Value = 8612.0625
$ws.Range("M136").Value = -8709.2001
$ws.Range("N136").Value = -13712.0625

$ws = $wb.Worksheets.Item("CUL")

# Row 140
$ws.Range("H140").Value = 1713.4286
$ws.Range("J140").Value = 2995
$ws.Range("L140").Value = 8985
$ws.Range("N140").Value = -19345

$ws = $wb.Worksheets.Item("GSM")

# Row 21
$ws.Range("H21").Value = 11665
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 14995
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 14995
$ws.Range("M21").Value = -9827
$ws.Range("N21").Value = -15341

# Row 30
$ws.Range("H30").Value = 11665
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 14995
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 14995
$ws.Range("M30").Value = -9895
$ws.Range("N30").Value = -15205

# Row 82
$ws.Range("H82").Value = 25298
$ws.Range("I82").Value = 25298
$ws.Range("K82").Value = 25298
$ws.Range("M82").Value = -24915

# Row 85
$ws.Range("H85").Value = 25298
$ws.Range("I85").Value = 25298
$ws.Range("K85").Value = 25298
$ws.Range("M85").Value = -23972

# Row 122
$ws.Range("H122").Value = 5268.2
$ws.Range("I122").Value = 5654.5713
$ws.Range("J122").Value = 4366.6665
$ws.Range("K122").Value = 16963.7139
$ws.Range("L122").Value = 13099.9995
$ws.Range("M122").Value = -14513.7139
$ws.Range("N122").Value = -17999.9995

# Row 126
$ws.Range("H126").Value = 2555.2104
$ws.Range("I126").Value = 2334
$ws.Range("K126").Value = 7002
$ws.Range("M126").Value = -4532

$ws = $wb.Worksheets.Item("LTW")

# Row 4
$ws.Range("H4").Value = 16000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 16000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 16000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -16226

# Row 25
$ws.Range("H25").Value = 6400
$ws.Range("J25").Value = 6400
$ws.Range("L25").Value = 6400
$ws.Range("N25").Value = -6860

# Row 28
$ws.Range("H28").Value = 16000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 16000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 16000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -16464

# Row 35
$ws.Range("H35").Value = 21213
$ws.Range("I35").Value = 5515
$ws.Range("J35").Value = 31678.334
$ws.Range("K35").Value = 5515
$ws.Range("L35").Value = 31678.334
$ws.Range("M35").Value = -5179
$ws.Range("N35").Value = -32350.334

# Row 37
$ws.Range("H37").Value = 16000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 16000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 16000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -16214

# Row 82
$ws.Range("H82").Value = 1620.4286
$ws.Range("I82").Value = 1708.3
$ws.Range("K82").Value = 1708.3
$ws.Range("M82").Value = -1347.3

# Row 85
$ws.Range("H85").Value = 1620.4286
$ws.Range("I85").Value = 1708.3
$ws.Range("K85").Value = 1708.3
$ws.Range("M85").Value = -460.3

$ws = $wb.Worksheets.Item("WVR")

# Row 2
$ws.Range("H2").Value = 30000
$ws.Range("I2").Value = 30000
$ws.Range("K2").Value = 30000
$ws.Range("M2").Value = -29888

# Row 4
$ws.Range("H4").Value = 805497.4
$ws.Range("I4").Value = 1335831
$ws.Range("J4").Value = 9997
$ws.Range("K4").Value = 1335831
$ws.Range("L4").Value = 9997
$ws.Range("M4").Value = -1335718
$ws.Range("N4").Value = -10223

# Row 62
$ws.Range("H62").Value = 9161.759
$ws.Range("J62").Value = 7738.1904
$ws.Range("L62").Value = 7738.1904
$ws.Range("N62").Value = -8986.190399999999

# Row 65
$ws.Range("H65").Value = 9161.759
$ws.Range("J65").Value = 7738.1904
$ws.Range("L65").Value = 38690.952
$ws.Range("N65").Value = -44930.952

# Row 107
$ws.Range("H107").Value = 2735.5
$ws.Range("I107").Value = 1748.5
$ws.Range("K107").Value = 5245.5
$ws.Range("M107").Value = -3325.5

# Row 126
$ws.Range("H126").Value = 2072.2856
$ws.Range("I126").Value = 1898.1765
$ws.Range("K126").Value = 5694.529500000001
$ws.Range("M126").Value = -3224.529500000001

